$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries to append after the existing data (rows 542-550).
# Columns: A=Date, B=Nom du joueur, C=Volume, D=Intensite, E=Fatigue,
#          F=Douleur, G=Localisation douleur, H=Plaisir, I=Charge (=C*D)
$rows = @(
    @{ Row=542; Date=45966; Name="Ilyes Boughanmi"; C=70; D=6; E=6; F=0; G=$null;     H=0 },
    @{ Row=543; Date=45966; Name="Amir Etien";       C=70; D=6; E=7; F=4; G="Mollet";  H=4 },
    @{ Row=544; Date=45966; Name="Yoan Zouma";       C=70; D=3; E=6; F=4; G="Cheville";H=5 },
    @{ Row=545; Date=45966; Name="Omar Benyounes";   C=70; D=5; E=6; F=0; G=$null;     H=5 },
    @{ Row=546; Date=45966; Name="Naim Ighbane";     C=70; D=5; E=3; F=7; G="Genou";   H=7 },
    @{ Row=547; Date=45966; Name="Karim Belmahi";    C=70; D=6; E=8; F=0; G=$null;     H=10 },
    @{ Row=548; Date=45966; Name="Malik Boussaid";   C=70; D=2; E=0; F=0; G=$null;     H=10 },
    @{ Row=549; Date=45966; Name="Naim Dhib";        C=70; D=6; E=5; F=2; G="Psoas ";  H=6 },
    @{ Row=550; Date=45966; Name="Mattheo Haon";     C=70; D=5; E=7; F=0; G=$null;     H=3 }
)

# Template row with text already present in every column (used to clone
# number formats / fonts) and a template empty-G cell to reuse for rows
# that have no "Localisation douleur" value.
$templateRow = 541
$emptyGTemplate = $ws.Range("G515")

foreach ($r in $rows) {
    $rowNum = $r.Row

    $src = $ws.Range("A" + $templateRow + ":I" + $templateRow)
    $dst = $ws.Range("A" + $rowNum + ":I" + $rowNum)
    $src.Copy($dst)

    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 2).Value = $r.Name
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Formula = "=C" + $rowNum + "*D" + $rowNum

    $gCell = $ws.Cells.Item($rowNum, 7)
    if ($r.G) {
        $gCell.Value = $r.G
    } else {
        $gCell.ClearContents()
        $emptyGTemplate.Copy($gCell)
    }
}

$excel.CutCopyMode = $false
$null = $ws.Range("L545").Select()
